$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the Cypher query text in B2:B4 and C2:C4 so that the WHERE clause
# filters on the lower-case literal 'unknown' instead of 'Unknown'.
$cells = @("B2", "C2", "B3", "C3", "B4", "C4")
foreach ($addr in $cells) {
    $range = $ws.Range($addr)
    $text = $range.Value2
    $newText = $text.Replace("g.platform in ['Unknown']", "g.platform in ['unknown']")
    $range.Value = $newText
}

# Select B5, matching the saved selection state in the workbook.
$ws.Range("B5").Select()
